# Removed spaces from headers
$wb = $excel.ActiveWorkbook
$wsSets = $wb.Worksheets.Item("Sets")
$wsParams = $wb.Worksheets.Item("Parameters")

# --- Parameters sheet header renames ---
# (ordering chosen to mirror how the shared-string table was rebuilt by the author)
$wsParams.Range("A1").Value2 = "Resident_ID"
$wsParams.Range("C1").Value2 = "Year_Level"

# "Week 1" .. "Week 52" live in columns O (15) through BN (66)
for ($w = 1; $w -le 52; $w++) {
    $col = 14 + $w
    $wsParams.Cells.Item(1, $col).Value2 = "Week_$w"
}

# --- Sets sheet header renames ---
$wsSets.Range("B1").Value2 = "Clinic_Group"
$wsSets.Range("D1").Value2 = "Number_of_Residents"

# Final Parameters header rename
$wsParams.Range("B1").Value2 = "Clinic_Groups"

# --- Selection state to match the saved workbook ---
[void]$wsSets.Range("D2").Select()
[void]$wsParams.Activate()
[void]$wsParams.Range("B1").Select()
